$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.645.87'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '2.296.53'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'300.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.57%  '
$ws.Range("D6").Value = "'95.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.81%  '
$ws.Range("E7").Value = '  -1.38%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -1.77%  '
$ws.Range("D10").Value = "'34.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.26%  '
$ws.Range("D11").Value = "'19.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.76%  '
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").Value = "'0.119"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D14").Value = "'6.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").Value = '2.656.65'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").Value = '2.302.19'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").Value = "'0.781"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").Value = '42.608.70'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").Value = "'12.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.85%  '
$ws.Range("E20").Value = '  -1.84%  '
$ws.Range("E21").Value = '  -0.75%  '
$ws.Range("D22").Value = "'67.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = "'235.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.85%  '
$ws.Range("D24").Value = "'2.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.37%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("E26").Value = '  -1.77%  '
$ws.Range("D27").Value = "'24.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.58%  '
$ws.Range("D28").Value = "'2.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +12.39%  '
$ws.Range("D29").Value = "'164.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.32%  '
$ws.Range("D30").Value = "'9.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("D31").Value = "'32.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.67%  '
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").Value = "'4.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.02%  '
$ws.Range("D34").Value = "'17.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.48%  '
$ws.Range("D35").Value = "'4.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.31%  '
$ws.Range("E36").Value = '  -2.26%  '
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").Value = "'0.0996"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.72%  '
$ws.Range("D39").Value = "'1.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("E41").Value = '  -1.32%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = "'19.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.76%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.959.66'
$ws.Range("E43").Value = '  -2.24%  '
$ws.Range("E44").Value = '  +5.03%  '
$ws.Range("E45").Value = '  -1.15%  '
$ws.Range("E46").Value = '  -3.50%  '
$ws.Range("D47").Value = "'2.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("D48").Value = "'2.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.82%  '
$ws.Range("D49").Value = '2.525.25'
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").Value = "'52.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.52%  '
$ws.Range("D51").Value = "'72.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.05%  '
